# Updates cryptos list prices / 1h volume percentages, and reorders a few
# coin rows, per the "Updated cryptos list" GitHub Actions commit.
#
# Several Price values (column D) are plain numeric-looking strings
# (e.g. "0.999", "1.00", "7.00"). Setting those via Range.Value directly
# would make Excel auto-convert them to numbers (losing the exact text,
# e.g. "1.00" -> 1), which does not match the source data (stored as
# text). To keep them as text without altering cell styles, we write a
# text formula that evaluates to the exact string, then replace the
# formula with its computed value in-place via Copy + PasteSpecial
# (xlPasteValues = -4163).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.179.95'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '2.228.58'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("E4").Value = '  -1.94%  '
$ws.Range("D5").Formula = '="299.14"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -2.27%  '
$ws.Range("D6").Formula = '="90.87"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  -3.41%  '
$ws.Range("D7").Formula = '="0.559"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  -1.91%  '
$ws.Range("D8").Formula = '="0.999"'
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  -0.43%  '
$ws.Range("D9").Formula = '="0.496"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  -4.65%  '
$ws.Range("D10").Formula = '="33.45"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  -3.51%  '
$ws.Range("D11").Formula = '="0.0781"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  -2.54%  '
$ws.Range("D12").Formula = '="7.00"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  -2.66%  '
$ws.Range("E13").Value = '  -0.56%  '
$ws.Range("D14").Value = '2.565.24'
$ws.Range("E14").Value = '  -0.57%  '
$ws.Range("D15").Value = '2.230.73'
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").Formula = '="13.41"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = '  -0.92%  '
$ws.Range("D17").Formula = '="0.780"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = '  -6.31%  '
$ws.Range("D18").Value = '44.010.86'
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").Formula = '="12.33"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  +3.47%  '
$ws.Range("D20").Value = '0.0₃0911'
$ws.Range("E20").Value = '  -4.31%  '
$ws.Range("D21").Formula = '="5.99"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -5.02%  '
$ws.Range("D22").Formula = '="64.21"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  -1.93%  '
$ws.Range("D23").Formula = '="236.00"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  -0.51%  '
$ws.Range("E24").Value = '  -4.60%  '
$ws.Range("D25").Formula = '="1.00"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("E26").Value = '  -6.88%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").Formula = '="2.27"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  +2.39%  '
$ws.Range("B28").Value = 'InjectiveProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D28").Formula = '="39.10"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  +1.82%  '
$ws.Range("D29").Formula = '="9.41"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  -3.83%  '
$ws.Range("D30").Formula = '="19.26"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  -3.40%  '
$ws.Range("D31").Formula = '="151.71"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  -0.90%  '
$ws.Range("D32").Formula = '="5.51"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  -7.23%  '
$ws.Range("D33").Formula = '="0.0766"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  -3.29%  '
$ws.Range("D34").Formula = '="2.51"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  -4.98%  '
$ws.Range("E35").Value = '  -1.32%  '
$ws.Range("E36").Value = '  -5.51%  '
$ws.Range("D37").Formula = '="2.87"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  -5.97%  '
$ws.Range("D38").Formula = '="1.69"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  -5.18%  '
$ws.Range("D39").Formula = '="0.0301"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  +0.51%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Formula = '="3.62"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  -3.42%  '
$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").Formula = '="3.17"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  -6.57%  '
$ws.Range("D42").Formula = '="13.61"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  -8.81%  '
$ws.Range("D43").Formula = '="0.999"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  -0.91%  '
$ws.Range("D44").Value = '1.806.10'
$ws.Range("E44").Value = '  +0.93%  '
$ws.Range("D45").Formula = '="1.80"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  +10.28%  '
$ws.Range("D46").Formula = '="0.186"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  -2.95%  '
$ws.Range("D47").Formula = '="68.50"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  -2.64%  '
$ws.Range("D48").Formula = '="94.75"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  -3.90%  '
$ws.Range("D49").Formula = '="73.43"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  -6.50%  '
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").Formula = '="4.63"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  -5.37%  '
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").Formula = '="7.76"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  -4.12%  '

$excel.CutCopyMode = 0
